$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '82.341.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.168.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '620.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.99%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.290'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +21.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.167.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.15%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.594'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000256'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.750.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.034.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.171.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.20%  '

$ws.Range("E19").Value = '  +10.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '437.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.330.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  -1.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '570.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.20%  '

$ws.Range("E35").Value = '  +22.30%  '

$ws.Range("E36").Value = '  -1.37%  '

$ws.Range("E37").Value = '  -1.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.24'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.404'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '160.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.25%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '185.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.765'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.17%  '
